$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump the "Date" value to the new publication timestamp.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value2 = "2024-03-22T16:25:12+00:00"

# ---------------------------------------------------------------------------
# 2. Elements sheet: columns AK (37) and AL (38) were swapped - both the
#    "Mapping: RIM Mapping" / "Mapping: Spécification métier vers
#    l'extension ROR CoordinateReliability" header columns, their widths,
#    and every data row underneath them.
# ---------------------------------------------------------------------------
$elem = $wb.Worksheets.Item("Elements")

# --- swap the column widths (AK <-> AL) -------------------------------------
# AK (37) was 24.98046875 wide / AL (38) was 76.828125 wide - after the edit
# AK takes AL's old (wider) width and AL takes AK's old (narrower) width.
$elem.Columns.Item(37).ColumnWidth = 76
$elem.Columns.Item(38).ColumnWidth = 24.166666666666668

# --- swap the cell contents row by row --------------------------------------
$lastRow = 6
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elem.Cells.Item($r, 37)
    $alCell = $elem.Cells.Item($r, 38)

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    if ($akVal -ne $alVal) {
        $akCell.Value2 = $alVal
        $alCell.Value2 = $akVal
    }
}
